# Apply "Update latest output (run 223)" changes to the workbook.
#
# Sheet "Schedule": update row 2 values and append a new row 3.
# Sheet "Detailed": update several cells in existing rows (9, 14-49) and
#                   append 48 new rows (50-97) for the second day's data.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Schedule")
$ws2 = $wb.Worksheets.Item("Detailed")

# ---------------------------------------------------------------------
# Sheet "Schedule"
# ---------------------------------------------------------------------

# Update existing row 2
$ws1.Cells.Item(2, 1).Value = 46060.16666666666
$ws1.Cells.Item(2, 2).Value = 46060.625
$ws1.Cells.Item(2, 5).Value = 1028.37523425
$ws1.Cells.Item(2, 6).Value = 24.73244911616162

# Add new row 3 - first copy the formatting from row 2 so the date-time
# number formats (columns A & B) carry over, then fill in the values.
$ws1.Range("A2:F2").Copy()
$ws1.Range("A3:F3").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws1.Cells.Item(3, 1).Value = 46061.29166666666
$ws1.Cells.Item(3, 2).Value = 46061.75
$ws1.Cells.Item(3, 3).Value = 11
$ws1.Cells.Item(3, 4).Value = 41.58
$ws1.Cells.Item(3, 5).Value = 1013.90383575
$ws1.Cells.Item(3, 6).Value = 24.38441163419914

# ---------------------------------------------------------------------
# Sheet "Detailed" - updates to existing rows
# ---------------------------------------------------------------------

$ws2.Cells.Item(9, 5).Value = "OFF"

$ws2.Cells.Item(14, 2).Value = 113.30932

$ws2.Cells.Item(15, 2).Value = 106.2694

$ws2.Cells.Item(16, 2).Value = 64.89
$ws2.Cells.Item(16, 3).Value = "historical"

$ws2.Cells.Item(17, 2).Value = 57.06007
$ws2.Cells.Item(17, 3).Value = "historical"

$ws2.Cells.Item(18, 2).Value = 36.05952
$ws2.Cells.Item(18, 3).Value = "historical"

$ws2.Cells.Item(19, 3).Value = "historical"

$ws2.Cells.Item(20, 3).Value = "historical"

$ws2.Cells.Item(21, 2).Value = 35.88
$ws2.Cells.Item(21, 3).Value = "historical"

$ws2.Cells.Item(22, 3).Value = "historical"

$ws2.Cells.Item(23, 2).Value = -8.52
$ws2.Cells.Item(23, 3).Value = "historical"

$ws2.Cells.Item(24, 2).Value = -5.51
$ws2.Cells.Item(24, 3).Value = "historical"

$ws2.Cells.Item(25, 2).Value = 0.51
$ws2.Cells.Item(25, 3).Value = "historical"

$ws2.Cells.Item(26, 2).Value = 3.86812
$ws2.Cells.Item(26, 3).Value = "historical"

$ws2.Cells.Item(27, 2).Value = 1.92032
$ws2.Cells.Item(27, 3).Value = "historical"

$ws2.Cells.Item(28, 3).Value = "historical"

$ws2.Cells.Item(29, 2).Value = 36.0601
$ws2.Cells.Item(29, 3).Value = "historical"

$ws2.Cells.Item(30, 3).Value = "historical"

$ws2.Cells.Item(31, 2).Value = 83.16965999999999
$ws2.Cells.Item(31, 3).Value = "historical"
$ws2.Cells.Item(31, 5).Value = "ON"

$ws2.Cells.Item(32, 2).Value = 95.50872
$ws2.Cells.Item(32, 3).Value = "historical"

$ws2.Cells.Item(33, 2).Value = 147.52

$ws2.Cells.Item(34, 2).Value = 147.52

$ws2.Cells.Item(35, 2).Value = 320.95786

$ws2.Cells.Item(36, 2).Value = 299.99

$ws2.Cells.Item(37, 2).Value = 169.19253

$ws2.Cells.Item(38, 2).Value = 171.60446

$ws2.Cells.Item(39, 2).Value = 189.62543

$ws2.Cells.Item(40, 2).Value = 147.21206

$ws2.Cells.Item(41, 2).Value = 137.88564

$ws2.Cells.Item(42, 2).Value = 137.56093

$ws2.Cells.Item(43, 2).Value = 138.41832

$ws2.Cells.Item(44, 2).Value = 155.00539

$ws2.Cells.Item(45, 2).Value = 137.62446

$ws2.Cells.Item(46, 2).Value = 138.42

$ws2.Cells.Item(47, 2).Value = 138.42

$ws2.Cells.Item(48, 2).Value = 139.72426

$ws2.Cells.Item(49, 2).Value = 138.42

# ---------------------------------------------------------------------
# Sheet "Detailed" - append new rows 50-97 (second day, 46061.*)
# ---------------------------------------------------------------------

# Copy the formatting of row 49 down across the new rows first, so that
# column A keeps its date-time number format and column D keeps its
# date number format.
$ws2.Range("A49:E49").Copy()
$ws2.Range("A50:E97").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$newRows = @(
    @(46061, 136.29152, "forecast", 46061, "OFF"),
    @(46061.02083333334, 108.89, "forecast", 46061, "OFF"),
    @(46061.04166666666, 138.25048, "forecast", 46061, "OFF"),
    @(46061.0625, 105.87452, "forecast", 46061, "OFF"),
    @(46061.08333333334, 105.79, "forecast", 46061, "OFF"),
    @(46061.10416666666, 105.35483, "forecast", 46061, "OFF"),
    @(46061.125, 105.79, "forecast", 46061, "OFF"),
    @(46061.14583333334, 103.50888, "forecast", 46061, "OFF"),
    @(46061.16666666666, 84.79000000000001, "forecast", 46061, "OFF"),
    @(46061.1875, 84.79000000000001, "forecast", 46061, "OFF"),
    @(46061.20833333334, 84.79000000000001, "forecast", 46061, "OFF"),
    @(46061.22916666666, 83.95312, "forecast", 46061, "OFF"),
    @(46061.25, 83.33698, "forecast", 46061, "OFF"),
    @(46061.27083333334, 78, "forecast", 46061, "OFF"),
    @(46061.29166666666, 60.78162, "forecast", 46061, "ON"),
    @(46061.3125, 58.76103, "forecast", 46061, "ON"),
    @(46061.33333333334, 41.60509, "forecast", 46061, "ON"),
    @(46061.35416666666, 57.08, "forecast", 46061, "ON"),
    @(46061.375, 56.97996, "forecast", 46061, "ON"),
    @(46061.39583333334, 56.97996, "forecast", 46061, "ON"),
    @(46061.41666666666, 36.07, "forecast", 46061, "ON"),
    @(46061.4375, 36.06035, "forecast", 46061, "ON"),
    @(46061.45833333334, 36.06004, "forecast", 46061, "ON"),
    @(46061.47916666666, 36.07, "forecast", 46061, "ON"),
    @(46061.5, 36.07, "forecast", 46061, "ON"),
    @(46061.52083333334, 40.54, "forecast", 46061, "ON"),
    @(46061.54166666666, 57.08, "forecast", 46061, "ON"),
    @(46061.5625, 45.56214, "forecast", 46061, "ON"),
    @(46061.58333333334, 36.07, "forecast", 46061, "ON"),
    @(46061.60416666666, 36.07, "forecast", 46061, "ON"),
    @(46061.625, 36.0601, "forecast", 46061, "ON"),
    @(46061.64583333334, 36.06036, "forecast", 46061, "ON"),
    @(46061.66666666666, 53.83207, "forecast", 46061, "ON"),
    @(46061.6875, 56.98, "forecast", 46061, "ON"),
    @(46061.70833333334, 61.98734, "forecast", 46061, "ON"),
    @(46061.72916666666, 67.14131, "forecast", 46061, "ON"),
    @(46061.75, 82.78167999999999, "forecast", 46061, "OFF"),
    @(46061.77083333334, 117.75187, "forecast", 46061, "OFF"),
    @(46061.79166666666, 121.09831, "forecast", 46061, "OFF"),
    @(46061.8125, 121.05307, "forecast", 46061, "OFF"),
    @(46061.83333333334, 144.6517, "forecast", 46061, "OFF"),
    @(46061.85416666666, 178.74691, "forecast", 46061, "OFF"),
    @(46061.875, 144.10367, "forecast", 46061, "OFF"),
    @(46061.89583333334, 138.42, "forecast", 46061, "OFF"),
    @(46061.91666666666, 108.89, "forecast", 46061, "OFF"),
    @(46061.9375, 107.7547, "forecast", 46061, "OFF"),
    @(46061.95833333334, 105.79, "forecast", 46061, "OFF"),
    @(46061.97916666666, 108.01, "forecast", 46061, "OFF")
)

$r = 50
foreach ($row in $newRows) {
    $ws2.Cells.Item($r, 1).Value = $row[0]
    $ws2.Cells.Item($r, 2).Value = $row[1]
    $ws2.Cells.Item($r, 3).Value = $row[2]
    $ws2.Cells.Item($r, 4).Value = $row[3]
    $ws2.Cells.Item($r, 5).Value = $row[4]
    $r = $r + 1
}
